$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the Hardware-ID value in row 3 to match the new licensed hardware ID
$ws.Range("A3").Value = "S36SNWAH859775Z"

# Move the active selection, mirroring the cursor position saved with the file
$ws.Range("A6").Select()
